$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(0, "model_1", "ModelType.LINEAR_REGRESSION", "ScalerType.Raw_", 0.3290896630758111, 1.230470725186108),
    @(2, "model_3", "ModelType.LINEAR_REGRESSION", "ScalerType.MinMaxScaler_", 0.3290896630758111, 1.230470725186108),
    @(1, "model_2", "ModelType.LINEAR_REGRESSION", "ScalerType.MaxAbsScaler_", 0.329089663075811, 1.230470725186108),
    @(5, "model_6", "ModelType.LINEAR_REGRESSION", "ScalerType.StandardScaler_", 0.3290896630758109, 1.230470725186108),
    @(4, "model_5", "ModelType.LINEAR_REGRESSION", "ScalerType.RobustScaler_", 0.3290896630758104, 1.230470725186109),
    @(9, "model_10", "ModelType.SVR", "ScalerType.Normalizer_", 0.3093889192968056, 1.266602510836622),
    @(7, "model_8", "ModelType.SVR", "ScalerType.MaxAbsScaler_", 0.2950808275515512, 1.292844002518694),
    @(8, "model_9", "ModelType.SVR", "ScalerType.MinMaxScaler_", 0.2868884014121623, 1.30786917052994),
    @(3, "model_4", "ModelType.LINEAR_REGRESSION", "ScalerType.Normalizer_", 0.2773614272766524, 1.325341941671152),
    @(11, "model_12", "ModelType.SVR", "ScalerType.StandardScaler_", 0.2771349193994743, 1.325757364264236),
    @(10, "model_11", "ModelType.SVR", "ScalerType.RobustScaler_", 0.2740432104482925, 1.331427656024333),
    @(24, "model_25", "ModelType.LIGHT_GBM", "ScalerType.Raw_", 0.2547040168462451, 1.366896347243373),
    @(25, "model_26", "ModelType.LIGHT_GBM", "ScalerType.MaxAbsScaler_", 0.2547040168462451, 1.366896347243373),
    @(26, "model_27", "ModelType.LIGHT_GBM", "ScalerType.MinMaxScaler_", 0.2547040168462451, 1.366896347243373),
    @(29, "model_30", "ModelType.LIGHT_GBM", "ScalerType.StandardScaler_", 0.2545472350692494, 1.367183889968859),
    @(28, "model_29", "ModelType.LIGHT_GBM", "ScalerType.RobustScaler_", 0.2529761818887426, 1.370065251068744),
    @(27, "model_28", "ModelType.LIGHT_GBM", "ScalerType.Normalizer_", 0.2117387969965522, 1.445695916003343),
    @(6, "model_7", "ModelType.SVR", "ScalerType.Raw_", 0.2068601848948719, 1.454643444518371),
    @(18, "model_19", "ModelType.RANDOM_FOREST", "ScalerType.Raw_", 0.0414956300501762, 1.75792725536659),
    @(22, "model_23", "ModelType.RANDOM_FOREST", "ScalerType.RobustScaler_", 0.02172418665476583, 1.794188706344167),
    @(23, "model_24", "ModelType.RANDOM_FOREST", "ScalerType.StandardScaler_", 0.01921383028087453, 1.7987927791358),
    @(20, "model_21", "ModelType.RANDOM_FOREST", "ScalerType.MinMaxScaler_", -0.004560913666249888, 1.842396409629628),
    @(19, "model_20", "ModelType.RANDOM_FOREST", "ScalerType.MaxAbsScaler_", -0.00814950784175661, 1.8489780045679),
    @(21, "model_22", "ModelType.RANDOM_FOREST", "ScalerType.Normalizer_", -0.01716134449921425, 1.865505997321742),
    @(33, "model_34", "ModelType.XGB", "ScalerType.Normalizer_", -0.02299239389664787, 1.876200325886644),
    @(30, "model_31", "ModelType.XGB", "ScalerType.Raw_", -0.1331679466264639, 2.078265765639371),
    @(31, "model_32", "ModelType.XGB", "ScalerType.MaxAbsScaler_", -0.1331679466264639, 2.078265765639371),
    @(32, "model_33", "ModelType.XGB", "ScalerType.MinMaxScaler_", -0.1331679466264639, 2.078265765639371),
    @(34, "model_35", "ModelType.XGB", "ScalerType.RobustScaler_", -0.1331679466264639, 2.078265765639371),
    @(35, "model_36", "ModelType.XGB", "ScalerType.StandardScaler_", -0.1331679466264639, 2.078265765639371),
    @(15, "model_16", "ModelType.DECISION_TREE", "ScalerType.Normalizer_", -0.5815337069533564, 2.900582716049383),
    @(14, "model_15", "ModelType.DECISION_TREE", "ScalerType.MinMaxScaler_", -0.6487801503206039, 3.023914814814815),
    @(12, "model_13", "ModelType.DECISION_TREE", "ScalerType.Raw_", -0.6670896755369637, 3.057495061728396),
    @(16, "model_17", "ModelType.DECISION_TREE", "ScalerType.RobustScaler_", -0.6670896755369637, 3.057495061728396),
    @(17, "model_18", "ModelType.DECISION_TREE", "ScalerType.StandardScaler_", -0.6688398507414688, 3.060704938271605),
    @(13, "model_14", "ModelType.DECISION_TREE", "ScalerType.MaxAbsScaler_", -0.6727440877361337, 3.067865432098765)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
